# This script updates the cryptocurrency price/volume snapshot data
# on the active worksheet to match the latest scrape, including two
# row re-orderings (PEPE / Binance-Peg BSC-USD and VeChain / RenderToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay plain text (Excel would otherwise
# auto-detect some of these strings, e.g. "543.43", as numbers). Forcing
# the number format to Text before the assignment keeps it a string, and
# resetting the style afterwards avoids leaving an extra format applied.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "61.630.66"
$ws.Range("E2").Value = "  -4.04%  "
$ws.Range("D3").Value = "2.976.14"
$ws.Range("E3").Value = "  -5.09%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws.Range("D5") "543.43"
$ws.Range("E5").Value = "  -4.73%  "
Set-TextValue $ws.Range("D6") "152.51"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("E7").Value = "  +0.10%  "
Set-TextValue $ws.Range("D8") "0.574"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "2.986.26"
$ws.Range("E9").Value = "  -5.11%  "
$ws.Range("E10").Value = "  -1.65%  "
Set-TextValue $ws.Range("D11") "6.16"
$ws.Range("E11").Value = "  -6.33%  "
Set-TextValue $ws.Range("D12") "0.372"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "3.497.83"
$ws.Range("E13").Value = "  -4.89%  "
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "61.699.05"
$ws.Range("E15").Value = "  -3.99%  "
Set-TextValue $ws.Range("D16") "23.73"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").Value = "2.975.85"
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("E18").Value = "  -3.76%  "
Set-TextValue $ws.Range("D19") "5.19"
$ws.Range("E19").Value = "  -0.40%  "
Set-TextValue $ws.Range("D20") "382.45"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("E22").Value = "  -5.66%  "
$ws.Range("E23").Value = "  -0.11%  "
Set-TextValue $ws.Range("D24") "65.74"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "3.101.35"
$ws.Range("E26").Value = "  -5.27%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "0.998"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0944"
$ws.Range("E29").Value = "  -5.66%  "
Set-TextValue $ws.Range("D30") "8.29"
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("E31").Value = "  -0.02%  "
Set-TextValue $ws.Range("D32") "1.73"
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("E33").Value = "  -2.77%  "
Set-TextValue $ws.Range("D34") "160.53"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  -1.38%  "
Set-TextValue $ws.Range("D36") "5.96"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("E39").Value = "  -6.06%  "
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").Value = "2.410.52"
$ws.Range("E42").Value = "  -9.29%  "
Set-TextValue $ws.Range("D43") "22.18"
$ws.Range("E43").Value = "  -5.67%  "
$ws.Range("E44").Value = "  -2.61%  "
Set-TextValue $ws.Range("D45") "0.0596"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "5.16"
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0249"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  +0.07%  "
Set-TextValue $ws.Range("D49") "270.46"
$ws.Range("E49").Value = "  -6.21%  "
Set-TextValue $ws.Range("D50") "19.86"
$ws.Range("E50").Value = "  -5.31%  "
Set-TextValue $ws.Range("D51") "0.0954"
$ws.Range("E51").Value = "  -1.67%  "
